$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.203.54"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.860.59"
$ws.Range("E3").Value = "  -0.77%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7132"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "240.58"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07741"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  -0.07%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "25.02"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08256"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.855.64"
$ws.Range("E12").Value = "  -1.28%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.236"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.7156"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "90.21"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "29.187.38"
$ws.Range("E16").Value = "  -1.54%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.872"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "244.40"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007811"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "2.109.52"
$ws.Range("E21").Value = "  -2.84%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.948"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1587"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "162.79"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.929"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.26"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -2.88%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.387"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.155"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05194"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.16%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.908"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.80%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.172"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.13%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7282"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.680"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01853"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.687"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "1.157.22"
$ws.Range("E40").Value = "  -1.88%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.9059"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E42").Value = "  +1.51%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "72.40"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "101.70"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "2.008.43"
$ws.Range("E46").Value = "  -2.36%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5241"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.768"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.312"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.871"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "7.063"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "
